# add function:makeclass save excel log
# ------------------------------------------------------------------
# This script reproduces a workbook save that happened after a new
# "makeclass" helper appended a merge-class row ("short") to the
# MergeData sheet, and bumped the related counters on ClassData.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsClass = $wb.Worksheets.Item("ClassData")
$wsMerge = $wb.Worksheets.Item("MergeData")
$wsName  = $wb.Worksheets.Item("NameData")

# --- ClassData (sheet1): bump merge counters / clear stale delete flags ---
$wsClass.Range("G14").Value = 0
$wsClass.Range("H14").Value = 0

$wsClass.Range("G32").Value = 0
$wsClass.Range("H32").Value = 0

$wsClass.Range("F43").Value = 4

$wsClass.Range("F45").Value = 4
$wsClass.Range("G45").Value = 0
$wsClass.Range("H45").Value = 0

$wsClass.Range("F50").Value = 5
$wsClass.Range("F51").Value = 5
$wsClass.Range("F52").Value = 5
$wsClass.Range("F53").Value = 5
$wsClass.Range("F54").Value = 5
$wsClass.Range("F55").Value = 5

$wsClass.Range("F60").Value = 6
$wsClass.Range("F61").Value = 7
$wsClass.Range("F62").Value = 6
$wsClass.Range("F63").Value = 7

$wsClass.Range("F65").Value = 8
$wsClass.Range("F66").Value = 8
$wsClass.Range("F67").Value = 8
$wsClass.Range("F69").Value = 8
$wsClass.Range("F70").Value = 8
$wsClass.Range("F73").Value = 8

# --- MergeData (sheet2): insert a new "short" merge-class row at row 5 ---
$wsMerge.Rows.Item(5).Insert()

# the new row inherits formatting from the row above it (row 4)
$wsMerge.Range("A4:C4").Copy()
$wsMerge.Range("A5:C5").PasteSpecial(-4122)

$wsMerge.Range("A5").Value = 4
$wsMerge.Range("B5").Value = 10
$wsMerge.Range("C5").Value = "short"

# re-number the idx column for the rows pushed down by the insert
$wsMerge.Range("A6").Value = 5
$wsMerge.Range("A7").Value = 6
$wsMerge.Range("A8").Value = 7
$wsMerge.Range("A9").Value = 8

# --- NameData (sheet3): move the saved cursor position ---
$wsName.Range("C16").Select() | Out-Null

# --- ClassData becomes the active sheet/tab again before saving ---
$wsClass.Activate() | Out-Null
